$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "21.092.52"
$ws.Range("E2").Value = "  +3.36%  "
$ws.Range("D3").Value = "1.538.65"
$ws.Range("E3").Value = "  +5.09%  "
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9682"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "282.09"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3631"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3194"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "40.72"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.29%  "
$ws.Range("E10").Value = "  +5.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06816"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.33%  "
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.683"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.73"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.367"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9678"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001043"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").Value = "1.536.87"
$ws.Range("E18").Value = "  +4.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06106"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.43"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.712"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.07"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.36"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.320"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.47%  "
$ws.Range("D25").Value = "21.162.34"
$ws.Range("E25").Value = "  +3.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "148.59"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.229"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +7.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.64"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.85%  "
$ws.Range("D29").Value = "1.705.92"
$ws.Range("E29").Value = "  +5.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.49"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.017"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8561"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +8.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.207"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08013"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.208"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.934"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05852"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.32%  "
$ws.Range("E39").Value = "  +3.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.66"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.739"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9680"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1919"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5451"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.54"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.568"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.45%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.11"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.47%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5441"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.870"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +6.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06567"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9921"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.04%  "
